{"js": "// The document currently starts with a table. The edit replaces that whole\n// table with two new paragraphs of placeholder body text (and otherwise\n// leaves the rest of the document - the trailing empty \"_GoBack\" paragraph\n// and the section properties - untouched).\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n\n  // Insert the two replacement paragraphs immediately before the table so\n  // they land at the very top of the document, in the same place the table\n  // used to be.\n  table.insertParagraph(\n    \"Video provides a powerful way to help you prove your point. When you click Online Video, you can paste in the embed code for the video you want to add. You can also type a keyword to search online for the video that best fits your document.\",\n    \"Before\"\n  );\n  table.insertParagraph(\n    \"To make your document look professionally produced, Word provides header, footer, cover page, and text box designs that complement each other. For example, you can add a matching cover page, header, and sidebar. Click Insert and then choose the elements you want from the different galleries.\",\n    \"Before\"\n  );\n\n  // Remove the now-superseded table entirely.\n  table.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The document currently opens with a single table. Replace it with two\n# paragraphs of placeholder body text inserted at the very top of the\n# document (where the table used to sit), then remove the table itself.\n$para2Text = \"To make your document look professionally produced, Word provides header, footer, cover page, and text box designs that complement each other. For example, you can add a matching cover page, header, and sidebar. Click Insert and then choose the elements you want from the different galleries.\"\n$para1Text = \"Video provides a powerful way to help you prove your point. When you click Online Video, you can paste in the embed code for the video you want to add. You can also type a keyword to search online for the video that best fits your document.\"\n\n$start = $d.Range(0, 0)\n$start.InsertBefore($para2Text + [char]13)\n\n$start = $d.Range(0, 0)\n$start.InsertBefore($para1Text + [char]13)\n\n$d.Tables(1).Delete()\n"}
